# Update the "取得日時" (retrieved datetime) timestamp stamped in column A
# for the data rows on the "ランサーズ" sheet.
#
# All rows that were fetched together share the same scrape timestamp
# "2025-09-09 12:38:05"; this run re-stamps them with the latest fetch
# time "2025-09-09 12:48:35" (commit: "Append: 2025-09-09 12:48 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-09 12:38:05"
$newTimestamp = "2025-09-09 12:48:35"

# Find the last used row based on column A (取得日時) and update every
# cell that still holds the previous scrape timestamp.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
